# don-quijote.docx edit
#
# The six verse/command lines are rewritten so the "git"-themed words
# (Repos, commit, log, head, master, push -> pull, remote, git commit -m "...")
# are wrapped in naive, escaped HTML markup (<p>, <em>, <code>), as if the
# paragraph text had been run through an HTML-escaping pass. Each sentence's
# final punctuation mark plus the opening "<" of its closing tag gets wrapped
# in its own run flanked by proofErr gramStart/gramEnd markers (mirroring
# what Word's grammar checker inserts around a run boundary straddling a
# comma/period). The stray "_GoBack" bookmark, which used to sit at the very
# start of the document, is moved to the end of the "pull flaco y remote
# corredor." paragraph.

$d = $word.ActiveDocument

# Run/paragraph formatting is identical everywhere in this document.
$rPr = '<w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="7030A0"/></w:rPr>'
$pPr = "<w:pPr>$rPr</w:pPr>"

function New-Run([string]$text, [bool]$preserve = $false) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $spaceAttr = ''
    if ($preserve) { $spaceAttr = ' xml:space="preserve"' }
    return "<w:r>$rPr<w:t$spaceAttr>$escaped</w:t></w:r>"
}

$gramStart = '<w:proofErr w:type="gramStart"/>'
$gramEnd = '<w:proofErr w:type="gramEnd"/>'
$goBack = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

# Paragraph 1: "En un lugar de los Repos," -> "<p>En un lugar de los <em>Repos</em>,</p>"
$p1 = (New-Run '<p>En un lugar de los <em>Repos</em') + $gramStart + (New-Run '>,<') + $gramEnd + (New-Run '/p>')

# Paragraph 2: "de cuyo commit no quiero acordarme," -> "<p>de cuyo <em>commit</em> no quiero acordarme,</p>"
$p2 = (New-Run '<p>de cuyo <em>commit</em> no quiero ' $true) + $gramStart + (New-Run 'acordarme,<') + $gramEnd + (New-Run '/p>')

# Paragraph 3: "no ha mucho log que vivía" -> "<p>no ha mucho <em>log</em> que vivía</p>" (single run)
$p3 = New-Run '<p>no ha mucho <em>log</em> que vivía</p>'

# Paragraph 4: "un hidalgo de los de head en master," -> "<p>un hidalgo de los de <em>head</em> en <em>master</em>,</p>"
$p4 = (New-Run '<p>un hidalgo de los de <em>head</em> en <em>master</em') + $gramStart + (New-Run '>,<') + $gramEnd + (New-Run '/p>')

# Paragraph 5: "push antiguo," -> "<p><em>push</em> antiguo,</p>"
$p5 = (New-Run '<p><em>push</em> ' $true) + $gramStart + (New-Run 'antiguo,<') + $gramEnd + (New-Run '/p>')

# Paragraph 6: "git flaco y remote corredor." -> "<p><em>pull</em> flaco y <em>remote</em> corredor.</p>"
# and the "_GoBack" bookmark now lives at the end of this paragraph.
$p6 = (New-Run '<p><em>pull</em> flaco y <em>remote</em> ' $true) + $gramStart + (New-Run 'corredor.<') + $gramEnd + (New-Run '/p>') + $goBack

# Paragraph 7: 'git commit -m "Don Quijote de la Mancha"' -> '<p><code>git commit -m &quot;Don Quijote de la Mancha&quot;</code> </p>'
# This paragraph loses its own <w:pPr> block entirely, and its old
# bookmarkEnd (the other half of the _GoBack bookmark) is gone, since the
# bookmark moved up into paragraph 6.
$p7 = (New-Run '<p><code>git commit -m &') + $gramStart + (New-Run 'quot;Don') + $gramEnd + (New-Run ' Quijote de la Mancha&quot;</code> </p>' $true)

$body  = "<w:p>$pPr$p1</w:p>"
$body += "<w:p>$pPr$p2</w:p>"
$body += "<w:p>$pPr$p3</w:p>"
$body += "<w:p>$pPr$p4</w:p>"
$body += "<w:p>$pPr$p5</w:p>"
$body += "<w:p>$pPr$p6</w:p>"
$body += "<w:p>$p7</w:p>"

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       "<w:body>$body</w:body>" +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole document body in one shot so that run/proofErr/bookmark
# placement comes out exactly as specified above (Find/Replace alone cannot
# precisely control proofErr insertion points or bookmark relocation).
$d.Content.InsertXML($xml)
